$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Escape-Xml($s) {
    $s = $s -replace '&', '&amp;'
    $s = $s -replace '<', '&lt;'
    $s = $s -replace '>', '&gt;'
    return $s
}

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Replace-ListBullet-Text finds the paragraph whose text equals $find exactly and
# rewrites it via InsertXML, preserving the paragraph's ListBullet pPr and the
# leading empty run that precedes the text run in the source document. (A plain
# Find/Replace on these paragraphs silently merges that leading empty run away,
# since the text run carries no distinguishing rPr.)
function Replace-ListBullet-Text($find, $replace) {
    foreach ($p in $d.Paragraphs) {
        $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($ptext -eq $find) {
            $rng = $d.Range($p.Range.Start, $p.Range.End)
            $safe = Escape-Xml $replace
            $xml = "<w:p $wns><w:pPr><w:pStyle w:val='ListBullet'/><w:spacing w:line='240' w:lineRule='auto'/><w:ind w:left='720'/></w:pPr><w:r><w:t>$safe</w:t></w:r></w:p>"
            $rng.InsertXML($xml)
        }
    }
}

# Heading + bold/italic byline paragraphs: the empty leading run survives a plain
# Find/Replace here because the text run carries distinguishing rPr (b/i), so a
# simple text substitution is enough.
Replace-Text "Play Book of Gods for Free - Review by Slot Expert" "Play Book of Gods Free: Review and Special Features"
Replace-Text "Read our review of Book of Gods - an Ancient Egyptian-themed online slot game. Play this visually amazing slot for free and explore its exciting features." "Read our review of Book of Gods, a visually stunning slot game. Play for free and enjoy special features."

# "What we like" / "What we don't like" ListBullet items: the text run has no
# distinguishing rPr, so a plain Find/Replace merges away the leading empty run.
# Rebuild these paragraphs via InsertXML instead, which keeps that run intact.
Replace-ListBullet-Text "Visually stunning graphics and great musical theme" "Outstanding visuals and sophisticated graphics"
Replace-ListBullet-Text "Flexible betting options cater to all bankrolls" "User-friendly control panel"
Replace-ListBullet-Text "Special features like Extra Scatter and Free Spins enhance gameplay" "Special features increase chances of winning"
Replace-ListBullet-Text "Competitive RTP value of 96.12%" "Flexible betting options for all bankrolls"
Replace-ListBullet-Text "Some may find the Egyptian theme overused in slot games" "Gamble feature can lead to potential loss of winnings"
